$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.607.20'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.98%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.819.05'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.47%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4672'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.46%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3588'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.15'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07120'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9012'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07778'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.831.20'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.248'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.327'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '87.42'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +3.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.008'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008545'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.651.22'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.19'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.012'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.61%  '
$ws.Range('E24').Value = '  +0.78%  '
$ws.Range('E25').Value = '  -2.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.17'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  +0.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.970'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.35'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.805'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08766'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.140'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.788'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +4.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7341'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.433'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('E36').Value = '  +1.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.075'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01930'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.909'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05107'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5069'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.795'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('E43').Value = '  -1.32%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.991'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4679'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.007'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.04'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '98.58'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.565'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06003'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.73'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.31%  '
